$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16687023639679
$ws.Range("B1").Value = 2.37345552444458
$ws.Range("D1").Value = 2.388326168060303
$ws.Range("E1").Value = 1.213736534118652
